$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined")

# Row 4: KNN / RAW
# (I4 previously held an empty placeholder cell with an explicit style;
# clear that leftover formatting before giving it a real value.)
$ws.Range("I4").ClearFormats()

$ws.Range("A4").Value = "Combined"
$ws.Range("B4").Value = 9589
$ws.Range("C4").Value = 2397
$ws.Range("D4").Value = 23433
$ws.Range("E4").Value = 120
$ws.Range("F4").Value = "KNN"
$ws.Range("G4").Value = "RAW"
$ws.Range("H4").Value = 0.7267
$ws.Range("I4").Value = 0.7267
$ws.Range("J4").Value = 0.4417
$ws.Range("K4").Value = 0.716

# Row 5: SVM / RAW
$ws.Range("A5").Value = "Combined"
$ws.Range("B5").Value = 9589
$ws.Range("C5").Value = 2397
$ws.Range("D5").Value = 23433
$ws.Range("E5").Value = 120
$ws.Range("F5").Value = "SVM"
$ws.Range("G5").Value = "RAW"
$ws.Range("H5").Value = 0.6742
$ws.Range("I5").Value = 0.6742
$ws.Range("J5").Value = 0.3044
$ws.Range("K5").Value = 0.612

# Row 6: RF / RAW
$ws.Range("A6").Value = "Combined"
$ws.Range("B6").Value = 9589
$ws.Range("C6").Value = 2397
$ws.Range("D6").Value = 23433
$ws.Range("E6").Value = 120
$ws.Range("F6").Value = "RF"
$ws.Range("G6").Value = "RAW"
$ws.Range("H6").Value = 0.3213
$ws.Range("I6").Value = 0.3213
$ws.Range("J6").Value = 0.04
$ws.Range("K6").Value = 0.1897

# Match the percentage number format used by the rest of column H
$ws.Range("H4:H6").NumberFormat = "0.00%"

# Column L stays empty on these new rows, but keeps a materialised cell
# (matching the rest of the table's per-row cell layout)
$ws.Range("L4:L6").Font.Bold = $false

# Update the selected cell shown in the sheet view
$ws.Activate()
$ws.Range("G16").Select()
